# The deck ships two theme parts:
#   - ppt/theme/theme1.xml  (bound to the slide master  -> currently the "Integral" palette)
#   - ppt/theme/theme2.xml  (bound to the notes master   -> currently the "Office Theme" palette)
#
# The authored change swaps the two colour palettes around: the slide
# master's theme becomes the stock "Office Theme" colours, and the notes
# master's theme becomes the "Integral" colours the slide master used to
# have. The font scheme and format scheme (fills/lines/effects/gradients)
# are byte-for-byte identical between the two themes already, so the
# swap is purely a change of the 12-slot colour scheme.
#
# PowerPoint's object model exposes the live theme colours through
# <Master>.Theme.ThemeColorScheme, whose 12 colour slots are read/write
# one-by-one via Colors(i).RGB (i = 1..12, in the fixed
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order). We drive the
# palette swap through that API.

$p = $ppt.ActivePresentation

# VBA-style RGB() helper: PowerPoint's *.RGB property (MsoRGBColor /
# RGBColor.RGB) packs the colour as R + G*256 + B*65536, i.e. byte-reversed
# from the "RRGGBB" hex notation used by <a:srgbClr val="RRGGBB"/>.
function RGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Stock PowerPoint "Office Theme" colours, in
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order.
$officeThemeColors = @(
    (RGB 0x00 0x00 0x00), # dk1      000000
    (RGB 0xFF 0xFF 0xFF), # lt1      FFFFFF
    (RGB 0x44 0x54 0x6A), # dk2      44546A
    (RGB 0xE7 0xE6 0xE6), # lt2      E7E6E6
    (RGB 0x5B 0x9B 0xD5), # accent1  5B9BD5
    (RGB 0xED 0x7D 0x31), # accent2  ED7D31
    (RGB 0xA5 0xA5 0xA5), # accent3  A5A5A5
    (RGB 0xFF 0xC0 0x00), # accent4  FFC000
    (RGB 0x44 0x72 0xC4), # accent5  4472C4
    (RGB 0x70 0xAD 0x47), # accent6  70AD47
    (RGB 0x05 0x63 0xC1), # hlink    0563C1
    (RGB 0x95 0x4F 0x72)  # folHlink 954F72
)

# Apply the new palette to the slide master's theme (ppt/theme/theme1.xml).
$masterColorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 0; $i -lt $officeThemeColors.Count; $i++) {
    $masterColorScheme.Colors($i + 1).RGB = $officeThemeColors[$i]
}
